$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6399.4
$ws.Range("J51").Value = 5499
$ws.Range("L51").Value = 5499
$ws.Range("N51").Value = -6467
$ws.Range("H86").Value = 637.25
$ws.Range("I86").Value = 637.25
$ws.Range("K86").Value = 637.25
$ws.Range("M86").Value = 485.75
$ws.Range("H89").Value = 637.25
$ws.Range("I89").Value = 637.25
$ws.Range("K89").Value = 3186.25
$ws.Range("M89").Value = 2429.75
$ws.Range("H132").Value = 6470.1
$ws.Range("I132").Value = 1837.625
$ws.Range("K132").Value = 5512.875
$ws.Range("M132").Value = -2982.875

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 36.333332
$ws.Range("I5").Value = 36.333332
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 36.333332
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 75.666668
$ws.Range("N5").ClearContents()
$ws.Range("H43").Value = 7550000
$ws.Range("I43").Value = 15000000
$ws.Range("K43").Value = 15000000
$ws.Range("M43").Value = -14999687
$ws.Range("H74").Value = 9992.25
$ws.Range("J74").Value = 9999
$ws.Range("L74").Value = 9999
$ws.Range("N74").Value = -11747
$ws.Range("H77").Value = 9992.25
$ws.Range("J77").Value = 9999
$ws.Range("L77").Value = 49995
$ws.Range("N77").Value = -58731
$ws.Range("H109").Value = 35000
$ws.Range("J109").Value = 35000
$ws.Range("L109").Value = 35000
$ws.Range("N109").Value = -37774

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 36.333332
$ws.Range("I4").Value = 36.333332
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 36.333332
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 78.666668
$ws.Range("N4").ClearContents()
$ws.Range("H86").Value = 4640.385
$ws.Range("I86").Value = 3053.5
$ws.Range("K86").Value = 3053.5
$ws.Range("M86").Value = -1930.5
$ws.Range("H89").Value = 4640.385
$ws.Range("I89").Value = 3053.5
$ws.Range("K89").Value = 15267.5
$ws.Range("M89").Value = -9651.5
$ws.Range("H123").Value = 20780
$ws.Range("J123").Value = 20780
$ws.Range("L123").Value = 20780
$ws.Range("N123").Value = -30580

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 7274
$ws.Range("I14").Value = 7274
$ws.Range("K14").Value = 7274
$ws.Range("M14").Value = -7104
$ws.Range("H15").Value = 1725
$ws.Range("J15").Value = 1725
$ws.Range("L15").Value = 1725
$ws.Range("N15").Value = -2065
$ws.Range("H22").Value = 1275.8235
$ws.Range("I22").Value = 751.1111
$ws.Range("K22").Value = 751.1111
$ws.Range("M22").Value = -401.1111
$ws.Range("H58").Value = 4694
$ws.Range("I58").Value = 1302.4
$ws.Range("J58").Value = 8933.5
$ws.Range("K58").Value = 1302.4
$ws.Range("L58").Value = 8933.5
$ws.Range("M58").Value = -1099.4
$ws.Range("N58").Value = -9339.5
$ws.Range("H131").Value = 38415.168
$ws.Range("J131").Value = 38415.168
$ws.Range("L131").Value = 38415.168
$ws.Range("N131").Value = -48495.168
$ws.Range("H136").Value = 4694
$ws.Range("I136").Value = 1302.4
$ws.Range("J136").Value = 8933.5
$ws.Range("K136").Value = 3907.2
$ws.Range("L136").Value = 26800.5
$ws.Range("M136").Value = -1357.2
$ws.Range("N136").Value = -31900.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 82
$ws.Range("I6").Value = 96.28570999999999
$ws.Range("J6").Value = 32
$ws.Range("K6").Value = 288.85713
$ws.Range("L6").Value = 96
$ws.Range("M6").Value = -175.85713
$ws.Range("N6").Value = -322
$ws.Range("H9").Value = 63.555557
$ws.Range("I9").Value = 54.25
$ws.Range("K9").Value = 162.75
$ws.Range("M9").Value = 61.25
$ws.Range("H19").Value = 3520.8
$ws.Range("I19").Value = 89
$ws.Range("J19").Value = 4378.75
$ws.Range("K19").Value = 267
$ws.Range("L19").Value = 13136.25
$ws.Range("M19").Value = -93
$ws.Range("N19").Value = -13484.25
$ws.Range("H60").Value = 1231.25
$ws.Range("I60").Value = 320
$ws.Range("J60").Value = 2750
$ws.Range("K60").Value = 960
$ws.Range("L60").Value = 8250
$ws.Range("M60").Value = -709
$ws.Range("N60").Value = -8752
$ws.Range("H98").Value = 199.5
$ws.Range("I98").Value = 122.666664
$ws.Range("J98").Value = 245.6
$ws.Range("K98").Value = 367.999992
$ws.Range("L98").Value = 736.8
$ws.Range("M98").Value = 1130.000008
$ws.Range("N98").Value = -3732.8
$ws.Range("H117").Value = 469.4
$ws.Range("J117").Value = 462
$ws.Range("L117").Value = 1386
$ws.Range("N117").Value = -8270
$ws.Range("H129").Value = 1262.4375
$ws.Range("I129").Value = 562.0769
$ws.Range("J129").Value = 4297.3335
$ws.Range("K129").Value = 1686.2307
$ws.Range("L129").Value = 12892.0005
$ws.Range("M129").Value = 3313.7693
$ws.Range("N129").Value = -22892.0005
$ws.Range("H134").Value = 3681.5454
$ws.Range("J134").Value = 2409
$ws.Range("L134").Value = 7227
$ws.Range("N134").Value = -17367

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 40
$ws.Range("J2").Value = 49.8
$ws.Range("L2").Value = 49.8
$ws.Range("N2").Value = -275.8
$ws.Range("H3").Value = 14389107
$ws.Range("J3").Value = 10000353
$ws.Range("L3").Value = 10000353
$ws.Range("N3").Value = -10000585
$ws.Range("H4").Value = 2712.5
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 3600
$ws.Range("K4").Value = 50
$ws.Range("L4").Value = 3600
$ws.Range("M4").Value = 62
$ws.Range("N4").Value = -3824
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H7").Value = 19167574
$ws.Range("J7").Value = 30001176
$ws.Range("L7").Value = 30001176
$ws.Range("N7").Value = -30001400
$ws.Range("H8").Value = 19167574
$ws.Range("J8").Value = 30001176
$ws.Range("L8").Value = 30001176
$ws.Range("N8").Value = -30001454
$ws.Range("H10").Value = 3100
$ws.Range("I10").Value = 2000
$ws.Range("J10").Value = 3833.3333
$ws.Range("K10").Value = 2000
$ws.Range("L10").Value = 3833.3333
$ws.Range("M10").Value = -1831
$ws.Range("N10").Value = -4171.3333
$ws.Range("H13").Value = 265.22223
$ws.Range("I13").Value = 139.5
$ws.Range("J13").Value = 301.14285
$ws.Range("K13").Value = 139.5
$ws.Range("L13").Value = 301.14285
$ws.Range("M13").Value = -0.5
$ws.Range("N13").Value = -579.14285
$ws.Range("H41").Value = 2285.25
$ws.Range("I41").Value = 2285.25
$ws.Range("K41").Value = 2285.25
$ws.Range("M41").Value = -1930.25
$ws.Range("H132").Value = 51078.816
$ws.Range("I132").Value = 64508.53
$ws.Range("J132").Value = 5417.8
$ws.Range("K132").Value = 193525.59
$ws.Range("L132").Value = 16253.4
$ws.Range("M132").Value = -190995.59
$ws.Range("N132").Value = -21313.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 6500
$ws.Range("I2").Value = 1000
$ws.Range("K2").Value = 1000
$ws.Range("M2").Value = -888
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H16").Value = 382
$ws.Range("I16").Value = 382
$ws.Range("K16").Value = 382
$ws.Range("M16").Value = -212
$ws.Range("H22").Value = 767.1429000000001
$ws.Range("I22").Value = 655
$ws.Range("J22").Value = 916.6667
$ws.Range("K22").Value = 655
$ws.Range("L22").Value = 916.6667
$ws.Range("M22").Value = -360
$ws.Range("N22").Value = -1506.6667
$ws.Range("H27").Value = 767.1429000000001
$ws.Range("I27").Value = 655
$ws.Range("J27").Value = 916.6667
$ws.Range("K27").Value = 655
$ws.Range("L27").Value = 916.6667
$ws.Range("M27").Value = -548
$ws.Range("N27").Value = -1130.6667
$ws.Range("H32").Value = 1644
$ws.Range("I32").Value = 1644
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1644
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1327
$ws.Range("N32").ClearContents()
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 141709.2
$ws.Range("I4").Value = 151828.08
$ws.Range("J4").Value = 45
$ws.Range("K4").Value = 151828.08
$ws.Range("L4").Value = 45
$ws.Range("M4").Value = -151715.08
$ws.Range("N4").Value = -271
$ws.Range("H15").Value = 9259
$ws.Range("J15").Value = 7775
$ws.Range("L15").Value = 7775
$ws.Range("N15").Value = -8351
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
